$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.158.95'
$ws.Range("E2").Value = '  +1.15%  '
$ws.Range("D3").Value = '3.561.54'
$ws.Range("E3").Value = '  +5.13%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '607.07'
$ws.Range("E5").Value = '  +2.06%  '
$ws.Range("D6").Value = '144.73'
$ws.Range("E6").Value = '  +2.53%  '
$ws.Range("D7").Value = '3.559.40'
$ws.Range("E7").Value = '  +5.10%  '
$ws.Range("E8").Value = '  +0.14%  '
$ws.Range("E9").Value = '  +4.38%  '
$ws.Range("E10").Value = '  +2.45%  '
$ws.Range("D11").Value = '7.97'
$ws.Range("E11").Value = '  +0.58%  '
$ws.Range("E12").Value = '  +1.62%  '
$ws.Range("D13").Value = '4.165.95'
$ws.Range("E13").Value = '  +5.20%  '
$ws.Range("D14").Value = '0.0000207'
$ws.Range("E14").Value = '  +4.25%  '
$ws.Range("D15").Value = '30.14'
$ws.Range("E15").Value = '  +1.97%  '
$ws.Range("D16").Value = '3.564.83'
$ws.Range("E16").Value = '  +5.46%  '
$ws.Range("D17").Value = '66.257.92'
$ws.Range("E17").Value = '  +1.34%  '
$ws.Range("E18").Value = '  -0.60%  '
$ws.Range("D19").Value = '11.41'
$ws.Range("E19").Value = '  +9.49%  '
$ws.Range("E20").Value = '  +1.61%  '
$ws.Range("D21").Value = '14.90'
$ws.Range("E21").Value = '  +2.01%  '
$ws.Range("D22").Value = '431.04'
$ws.Range("E22").Value = '  +4.16%  '
$ws.Range("D23").Value = '0.610'
$ws.Range("E23").Value = '  +5.58%  '
$ws.Range("E24").Value = '  +2.55%  '
$ws.Range("D25").Value = '3.703.63'
$ws.Range("E25").Value = '  +5.21%  '
$ws.Range("E26").Value = '  -0.08%  '
$ws.Range("E27").Value = '  +8.84%  '
$ws.Range("D28").Value = '2.53'
$ws.Range("E28").Value = '  +4.81%  '
$ws.Range("D29").Value = '8.01'
$ws.Range("E29").Value = '  +3.26%  '
$ws.Range("E30").Value = '  -1.05%  '
$ws.Range("E31").Value = '  -0.17%  '
$ws.Range("E32").Value = '  +1.60%  '
$ws.Range("D33").Value = '25.52'
$ws.Range("E33").Value = '  +4.91%  '
$ws.Range("B34").Value = 'RenzoRestakedETH'
$ws.Range("C34").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D34").Value = '3.556.16'
$ws.Range("E34").Value = '  +5.13%  '
$ws.Range("B35").Value = 'Kaspa'
$ws.Range("C35").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D35").Value = '0.155'
$ws.Range("E35").Value = '  -2.68%  '
$ws.Range("E36").Value = '  +0.04%  '
$ws.Range("D37").Value = '1.75'
$ws.Range("E37").Value = '  +4.51%  '
$ws.Range("D38").Value = '7.90'
$ws.Range("E38").Value = '  +5.30%  '
$ws.Range("D39").Value = '5.62'
$ws.Range("E39").Value = '  +1.49%  '
$ws.Range("E40").Value = '  +0.08%  '
$ws.Range("D41").Value = '170.13'
$ws.Range("E41").Value = '  +1.48%  '
$ws.Range("E42").Value = '  +0.17%  '
$ws.Range("E43").Value = '  +3.86%  '
$ws.Range("E44").Value = '  +3.21%  '
$ws.Range("D45").Value = '1.93'
$ws.Range("E45").Value = '  +1.43%  '
$ws.Range("D46").Value = '46.18'
$ws.Range("E46").Value = '  +1.97%  '
$ws.Range("E47").Value = '  +3.10%  '
$ws.Range("E48").Value = '  -2.95%  '
$ws.Range("E49").Value = '  +5.35%  '
$ws.Range("E50").Value = '  +1.76%  '
$ws.Range("D51").Value = '23.38'
$ws.Range("E51").Value = '  +16.31%  '
